$d = $word.ActiveDocument

# 1) Apply strikethrough to the "Using create-react-app, create a new React project." paragraph
#    (including the paragraph mark itself, matching Word's behavior when the whole line incl.
#    pilcrow is selected and Strikethrough is toggled on).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Using create-react-app*") {
        $p.Range.Font.StrikeThrough = 1
    }
}

# 2) Replace the split "don't" run sequence (with proofErr wrappers) by retyping the sentence
#    as a single contiguous run, matching Word's behavior after an edit that replaces the
#    selected text with freshly-typed text that no longer carries per-run proofErr markers.
$apostrophe = [char]0x2019
$oldText = "a Navigation component that contains links styled like a navbar. The links don" + $apostrophe + "t have to go anywhere"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)
